$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4935.636
$ws.Range("I40").Value = 6618.857
$ws.Range("K40").Value = 6618.857
$ws.Range("M40").Value = -6443.857
$ws.Range("H70").Value = 8274.519
$ws.Range("I70").Value = 1629.2
$ws.Range("J70").Value = 9784.817999999999
$ws.Range("K70").Value = 4887.6
$ws.Range("L70").Value = 29354.454
$ws.Range("M70").Value = -4617.6
$ws.Range("N70").Value = -29894.454
$ws.Range("H73").Value = 8274.519
$ws.Range("I73").Value = 1629.2
$ws.Range("J73").Value = 9784.817999999999
$ws.Range("K73").Value = 4887.6
$ws.Range("L73").Value = 29354.454
$ws.Range("M73").Value = -3951.6
$ws.Range("N73").Value = -31226.454
$ws.Range("H88").Value = 1170.6471
$ws.Range("I88").Value = 1083.4546
$ws.Range("J88").Value = 1330.5
$ws.Range("K88").Value = 1083.4546
$ws.Range("L88").Value = 1330.5
$ws.Range("M88").Value = -677.4546
$ws.Range("N88").Value = -2142.5
$ws.Range("H91").Value = 1170.6471
$ws.Range("I91").Value = 1083.4546
$ws.Range("J91").Value = 1330.5
$ws.Range("K91").Value = 1083.4546
$ws.Range("L91").Value = 1330.5
$ws.Range("M91").Value = 320.5454
$ws.Range("N91").Value = -4138.5
$ws.Range("H92").Value = 3027.639
$ws.Range("I92").Value = 2835.963
$ws.Range("J92").Value = 3602.6667
$ws.Range("K92").Value = 2835.963
$ws.Range("L92").Value = 3602.6667
$ws.Range("M92").Value = -1587.963
$ws.Range("N92").Value = -6098.6667
$ws.Range("H96").Value = 2455.25
$ws.Range("I96").Value = 2612
$ws.Range("K96").Value = 7836
$ws.Range("M96").Value = -6463
$ws.Range("H99").Value = 172619860
$ws.Range("I99").Value = 11905549
$ws.Range("K99").Value = 35716647
$ws.Range("M99").Value = -35715149
$ws.Range("H135").Value = 930.5227
$ws.Range("I135").Value = 602.325
$ws.Range("K135").Value = 5420.925
$ws.Range("M135").Value = -2885.925
$ws.Range("H137").Value = 38458.418
$ws.Range("I137").Value = 57419.1
$ws.Range("J137").Value = 3984.4546
$ws.Range("K137").Value = 172257.3
$ws.Range("L137").Value = 11953.3638
$ws.Range("M137").Value = -169707.3
$ws.Range("N137").Value = -17053.3638

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 73988.78999999999
$ws.Range("I74").Value = 101222.85
$ws.Range("K74").Value = 101222.85
$ws.Range("M74").Value = -100348.85
$ws.Range("H75").Value = 40000
$ws.Range("J75").Value = 40000
$ws.Range("L75").Value = 40000
$ws.Range("N75").Value = -41748
$ws.Range("H77").Value = 73988.78999999999
$ws.Range("I77").Value = 101222.85
$ws.Range("K77").Value = 506114.25
$ws.Range("M77").Value = -501746.25
$ws.Range("H78").Value = 40000
$ws.Range("J78").Value = 40000
$ws.Range("L78").Value = 120000
$ws.Range("N78").Value = -128736
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 2513.513
$ws.Range("J132").Value = 2870.75
$ws.Range("L132").Value = 8612.25
$ws.Range("N132").Value = -13672.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1550.3
$ws.Range("I58").Value = 1603.3513
$ws.Range("K58").Value = 1603.3513
$ws.Range("M58").Value = -1400.3513
$ws.Range("H132").Value = 2812.7737
$ws.Range("I132").Value = 2937.0625
$ws.Range("K132").Value = 8811.1875
$ws.Range("M132").Value = -6281.1875
$ws.Range("H134").Value = 7345.4116
$ws.Range("I134").Value = 8423.143
$ws.Range("K134").Value = 25269.429
$ws.Range("M134").Value = -22734.429
$ws.Range("H136").Value = 1550.3
$ws.Range("I136").Value = 1603.3513
$ws.Range("K136").Value = 4810.0539
$ws.Range("M136").Value = -2260.0539

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 5555723
$ws.Range("I2").Value = 257.1111
$ws.Range("J2").Value = 11111189
$ws.Range("K2").Value = 1542.6666
$ws.Range("L2").Value = 66667134
$ws.Range("M2").Value = -1429.6666
$ws.Range("N2").Value = -66667360
$ws.Range("H113").Value = 539.4
$ws.Range("J113").Value = 481.33334
$ws.Range("L113").Value = 1444.00002
$ws.Range("N113").Value = -5784.000019999999
$ws.Range("H122").Value = 1434.1428
$ws.Range("J122").Value = 925
$ws.Range("L122").Value = 8325
$ws.Range("N122").Value = -13225

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 10675.125
$ws.Range("I43").Value = 7574.3335
$ws.Range("J43").Value = 19977.5
$ws.Range("K43").Value = 7574.3335
$ws.Range("L43").Value = 19977.5
$ws.Range("M43").Value = -7423.3335
$ws.Range("N43").Value = -20279.5
$ws.Range("H80").Value = 3582.077
$ws.Range("J80").Value = 4088.2942
$ws.Range("L80").Value = 4088.2942
$ws.Range("N80").Value = -6084.2942
$ws.Range("H83").Value = 3582.077
$ws.Range("J83").Value = 4088.2942
$ws.Range("L83").Value = 20441.471
$ws.Range("N83").Value = -30425.471
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H122").Value = 3848577.5
$ws.Range("I122").Value = 6251494
$ws.Range("K122").Value = 18754482
$ws.Range("M122").Value = -18752032
$ws.Range("H132").Value = 44777.81
$ws.Range("I132").Value = 48263.457
$ws.Range("J132").Value = 2950
$ws.Range("K132").Value = 144790.371
$ws.Range("L132").Value = 8850
$ws.Range("M132").Value = -142260.371
$ws.Range("N132").Value = -13910

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8519.091
$ws.Range("I7").Value = 8776.25
$ws.Range("J7").Value = 5947.5
$ws.Range("K7").Value = 8776.25
$ws.Range("L7").Value = 5947.5
$ws.Range("M7").Value = -8664.25
$ws.Range("N7").Value = -6171.5
$ws.Range("H40").Value = 16294.1875
$ws.Range("J40").Value = 4443.4
$ws.Range("L40").Value = 4443.4
$ws.Range("N40").Value = -4715.4
$ws.Range("H126").Value = 8519.091
$ws.Range("I126").Value = 8776.25
$ws.Range("J126").Value = 5947.5
$ws.Range("K126").Value = 26328.75
$ws.Range("L126").Value = 17842.5
$ws.Range("M126").Value = -23858.75
$ws.Range("N126").Value = -22782.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 14268
$ws.Range("J41").Value = 14268
$ws.Range("L41").Value = 14268
$ws.Range("N41").Value = -15048
$ws.Range("H81").Value = 2529.6
$ws.Range("I81").Value = 2588.4443
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 5176.8886
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -4115.8886
$ws.Range("N81").Value = -6122
$ws.Range("H84").Value = 2529.6
$ws.Range("I84").Value = 2588.4443
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 25884.443
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -20580.443
$ws.Range("N84").Value = -30608
$ws.Range("H136").Value = 199343
$ws.Range("I136").Value = 230419.95
$ws.Range("K136").Value = 691259.8500000001
$ws.Range("M136").Value = -688709.8500000001

